$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before V (shifts old V,W,X -> W,X,Y) to make room for
# the new "TakeBets" payout column.
$ws.Columns("V").Insert()

# Give the freshly inserted column the narrower width used for the new
# TakeBets column (stored column width "9" once Excel applies its internal
# padding of ~0.71 characters).
$ws.Columns("V").ColumnWidth = 8.14

# New "TakeBets" header/value pair in the freshly inserted column.
$ws.Range("V1").Value = "TakeBets"
$ws.Range("V2").Value = "B3;B5"

# The old payAmt1 header text drops the trailing ";Player" suffix; it now
# lives (after the column insert) in W1.
$ws.Range("W1").Value = "payAmt1;Antenna"

# Fill in the new payout amounts that go with the shifted payAmt columns.
$ws.Range("W2").Value = "P1;100"
$ws.Range("X2").Value = "P2;100"
$ws.Range("Y2").Value = "P6;100"

# The sample card values in row 2 were reshuffled.
$ws.Range("P2").Value = "4d"
$ws.Range("Q2").Value = "2s"
$ws.Range("R2").Value = "4d"
$ws.Range("S2").Value = "3s"

# Update the view/selection to match the new layout.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 14
$ws.Range("W1:Y2").Select()
